$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dico")

# --- Table body edits -------------------------------------------------
# The big rich-text "question client" note that used to live in the merged
# G2:J6 cell is being moved into a floating text box instead, so the cell
# itself goes blank (its fill/alignment style "1" is left untouched).
$ws.Range("G2").Value = ""

# The "departement" block gains a constraint entry that was missing before
# (D7 = "AN(50)", matching the type already used by D5/D6).
$ws.Range("D7").Value = "AN(50)"

# Excel explicitly recorded a custom (but unchanged) row height for row 7
# once the sheet was touched again.
$ws.Rows.Item(7).RowHeight = 14.4

# The author's cursor ended up on C16 when the file was last saved.
$ws.Range("C16").Select()

# --- Floating text box with the relocated question/answer text --------
$shp = $ws.Shapes.AddTextbox(1, 512.4, 18.0, 193.8, 157.8)
$shp.Name = "ZoneTexte 1"
$shp.Fill.ForeColor.RGB = 16777215

$question = "es-que une ville peu avoir plusieurs plage? " + `
  'Reponse client "oui". ' + `
  "Combien de touristes en moyen avez-vous dans les villes? " + `
  "Quel est la longueur maximum d'une de vos plages? " + `
  "Comment nomm" + [char]0xE9 + " vous vos plages?"

$shp.TextFrame.Characters().Text = $question
